$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

# D12: new value 1, formatted as percentage (100%)
$ws.Range("D12").Value = 1
$ws.Range("D12").NumberFormat = "0%"

# E12: new note text
$ws.Range("E12").Value = "100%(17/06/2010)"

# Update the visible selection/scroll state to match the saved view
$ws.Range("G14").Select() | Out-Null
